# Auto-generated Excel COM-interop script
# Applies the numeric cell updates described in the commit diff
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 125000180
$ws.Range("J9").Value = 249.5
$ws.Range("L9").Value = 249.5
$ws.Range("N9").Value = -587.5
$ws.Range("H12").Value = 665.9167
$ws.Range("J12").Value = 850
$ws.Range("L12").Value = 850
$ws.Range("N12").Value = -1190
$ws.Range("H40").Value = 2574.0476
$ws.Range("I40").Value = 2332.6428
$ws.Range("J40").Value = 3056.8572
$ws.Range("K40").Value = 2332.6428
$ws.Range("L40").Value = 3056.8572
$ws.Range("M40").Value = -2157.6428
$ws.Range("N40").Value = -3406.8572
$ws.Range("H69").Value = 33250
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 33250
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 99750
$ws.Range("N69").Value = -101498
$ws.Range("H72").Value = 33250
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 33250
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 299250
$ws.Range("N72").Value = -307986
$ws.Range("H137").Value = 874106
$ws.Range("I137").Value = 3438959.2
$ws.Range("J137").Value = 19154.889
$ws.Range("K137").Value = 10316877.6
$ws.Range("L137").Value = 57464.667
$ws.Range("M137").Value = -10314327.6
$ws.Range("N137").Value = -62564.667
$ws.Range("H138").Value = 5376.3623
$ws.Range("I138").Value = 1297
$ws.Range("J138").Value = 6554.844
$ws.Range("K138").Value = 3891
$ws.Range("L138").Value = 19664.532
$ws.Range("M138").Value = 1249
$ws.Range("N138").Value = -29944.532
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2983.0625
$ws.Range("I74").Value = 2964.1538
$ws.Range("J74").Value = 3065
$ws.Range("K74").Value = 2964.1538
$ws.Range("L74").Value = 3065
$ws.Range("M74").Value = -2090.1538
$ws.Range("N74").Value = -4813
$ws.Range("H77").Value = 2983.0625
$ws.Range("I77").Value = 2964.1538
$ws.Range("J77").Value = 3065
$ws.Range("K77").Value = 14820.769
$ws.Range("L77").Value = 15325
$ws.Range("M77").Value = -10452.769
$ws.Range("N77").Value = -24061
$ws.Range("H102").Value = 7348.793
$ws.Range("I102").Value = 5258.2144
$ws.Range("K102").Value = 5258.2144
$ws.Range("M102").Value = -3636.2144
$ws.Range("H132").Value = 2559.5
$ws.Range("I132").Value = 2144.5715
$ws.Range("J132").Value = 4302.2
$ws.Range("K132").Value = 6433.7145
$ws.Range("L132").Value = 12906.6
$ws.Range("M132").Value = -3903.7145
$ws.Range("N132").Value = -17966.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 23256.875
$ws.Range("I99").Value = 25900.857
$ws.Range("K99").Value = 25900.857
$ws.Range("M99").Value = -24402.857
$ws.Range("H105").Value = 58182.555
$ws.Range("I105").Value = 101934.7
$ws.Range("J105").Value = 3492.375
$ws.Range("K105").Value = 101934.7
$ws.Range("L105").Value = 3492.375
$ws.Range("M105").Value = -100187.7
$ws.Range("N105").Value = -6986.375
$ws.Range("H134").Value = 6388.387
$ws.Range("I134").Value = 7739.8
$ws.Range("K134").Value = 23219.4
$ws.Range("M134").Value = -20684.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 103070.4
$ws.Range("I16").Value = 3315
$ws.Range("K16").Value = 3315
$ws.Range("M16").Value = -3028
$ws.Range("H113").Value = 103070.4
$ws.Range("I113").Value = 3315
$ws.Range("K113").Value = 3315
$ws.Range("M113").Value = -1145
$ws.Range("H132").Value = 40878.73
$ws.Range("I132").Value = 11898.9
$ws.Range("J132").Value = 137478.17
$ws.Range("K132").Value = 35696.7
$ws.Range("L132").Value = 412434.51
$ws.Range("M132").Value = -33166.7
$ws.Range("N132").Value = -417494.51
$ws.Range("H133").Value = 35000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 320.2
$ws.Range("I14").Value = 320.2
$ws.Range("K14").Value = 960.5999999999999
$ws.Range("M14").Value = -787.5999999999999
$ws.Range("H137").Value = 8452.5
$ws.Range("J137").Value = 10767.923
$ws.Range("L137").Value = 32303.769
$ws.Range("N137").Value = -42503.769
$ws.Range("H141").Value = 1000
$ws.Range("I141").Value = 1000
$ws.Range("K141").Value = 3000
$ws.Range("M141").Value = 2180
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 10609.192
$ws.Range("I122").Value = 11503.565
$ws.Range("K122").Value = 34510.695
$ws.Range("M122").Value = -32060.695
$ws.Range("H132").Value = 4384.107
$ws.Range("J132").Value = 6923.8335
$ws.Range("L132").Value = 20771.5005
$ws.Range("N132").Value = -25831.5005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2116.4375
$ws.Range("I22").Value = 2648.7144
$ws.Range("K22").Value = 2648.7144
$ws.Range("M22").Value = -2353.7144
$ws.Range("H27").Value = 2116.4375
$ws.Range("I27").Value = 2648.7144
$ws.Range("K27").Value = 2648.7144
$ws.Range("M27").Value = -2541.7144
$ws.Range("H100").Value = 5070.5293
$ws.Range("J100").Value = 3500
$ws.Range("L100").Value = 3500
$ws.Range("N100").Value = -4582
$ws.Range("H136").Value = 9339.1875
$ws.Range("I136").Value = 2361
$ws.Range("J136").Value = 11665.25
$ws.Range("K136").Value = 7083
$ws.Range("L136").Value = 34995.75
$ws.Range("M136").Value = -4533
$ws.Range("N136").Value = -40095.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("H62").Value = 110749.84
$ws.Range("I62").Value = 301861.22
$ws.Range("J62").Value = 3249.6875
$ws.Range("K62").Value = 301861.22
$ws.Range("L62").Value = 3249.6875
$ws.Range("M62").Value = -301237.22
$ws.Range("N62").Value = -4497.6875
$ws.Range("H65").Value = 110749.84
$ws.Range("I65").Value = 301861.22
$ws.Range("J65").Value = 3249.6875
$ws.Range("K65").Value = 1509306.1
$ws.Range("L65").Value = 16248.4375
$ws.Range("M65").Value = -1506186.1
$ws.Range("N65").Value = -22488.4375
$ws.Range("H122").Value = 5743.4165
$ws.Range("I122").Value = 4436.3335
$ws.Range("K122").Value = 13309.0005
$ws.Range("M122").Value = -10859.0005
$ws.Range("H132").Value = 18431.258
$ws.Range("I132").Value = 21623.46
$ws.Range("K132").Value = 64870.38
$ws.Range("M132").Value = -62340.38
$ws.Range("H136").Value = 4316.4736
$ws.Range("I136").Value = 2693.6667
$ws.Range("K136").Value = 8081.000100000001
$ws.Range("M136").Value = -5531.000100000001

# Cells that are fully removed (no cached value remains) in the target state
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M69").ClearContents()
$ws.Range("M72").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M21").ClearContents()
$ws.Range("M35").ClearContents()

Write-Host "Applied all Sheets updates"
